$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.649812734082397
$ws1.Range("C2").Value = 0.613314447592068
$ws1.Range("D2").Value = 0.8108614232209738
$ws1.Range("E2").Value = 0.6983870967741935
$ws1.Range("F2").Value = 0.7617874736101337
$ws1.Range("G2").Value = 0.8009391007398976
$ws1.Range("H2").Value = 0.6820512280997069
$ws1.Range("I2").Value = 433
$ws1.Range("J2").Value = 273
$ws1.Range("K2").Value = 261
$ws1.Range("L2").Value = 101

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.7209944751381215
$ws2.Range("C2").Value = 0.4887640449438202
$ws2.Range("D2").Value = 0.5825892857142857

$ws2.Range("B3").Value = 0.613314447592068
$ws2.Range("C3").Value = 0.8108614232209738
$ws2.Range("D3").Value = 0.6983870967741935

$ws2.Range("B4").Value = 0.649812734082397
$ws2.Range("C4").Value = 0.649812734082397
$ws2.Range("D4").Value = 0.649812734082397
$ws2.Range("E4").Value = 0.649812734082397

$ws2.Range("B5").Value = 0.6671544613650948
$ws2.Range("C5").Value = 0.6498127340823969
$ws2.Range("D5").Value = 0.6404881912442396

$ws2.Range("B6").Value = 0.6671544613650948
$ws2.Range("C6").Value = 0.649812734082397
$ws2.Range("D6").Value = 0.6404881912442396

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 261
$ws3.Range("C2").Value = 273

$ws3.Range("B3").Value = 101
$ws3.Range("C3").Value = 433
